$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data refresh effectively re-maps each data row (2-11) to the
# values that used to belong to a different row, while columns A, B, C, E, F,
# G, H, I, J, K, R stay the same for every row. Only D (Fecha), L (Calidad),
# M (Volumen), N/O/P (precios), Q (Unidad), S (Precio $/Kg) and T (Kg/unidad)
# move around, following this row permutation (new row <- old row):
#   2<-6  3<-7  4<-9  5<-10  6<-5  7<-8  8<-2  9<-11  10<-3  11<-4

$rows = @{
    2  = @{ D = 44400; L = "Primera"; M = 25; N = 1500;  O = 1500;  P = 1500;  Q = "`$/kilo (en caja de 15 kilos)"; S = 1500; T = 1 }
    3  = @{ D = 44309; L = "Primera"; M = 10; N = 1600;  O = 1600;  P = 1600;  Q = "`$/kilo (en caja de 15 kilos)"; S = 1600; T = 1 }
    4  = @{ D = 44371; L = "Primera"; M = 20; N = 1800;  O = 1800;  P = 1800;  Q = "`$/kilo (en caja de 15 kilos)"; S = 1800; T = 1 }
    5  = @{ D = 44371; L = "Segunda"; M = 30; N = 1200;  O = 1200;  P = 1200;  Q = "`$/kilo (en caja de 15 kilos)"; S = 1200; T = 1 }
    6  = @{ D = 44336; L = "Primera"; M = 10; N = 1500;  O = 1500;  P = 1500;  Q = "`$/kilo (en caja de 15 kilos)"; S = 1500; T = 1 }
    7  = @{ D = 44343; L = "Primera"; M = 20; N = 1700;  O = 1700;  P = 1700;  Q = "`$/kilo (en caja de 15 kilos)"; S = 1700; T = 1 }
    8  = @{ D = 44195; L = "Primera"; M = 20; N = 15000; O = 15000; P = 15000; Q = "`$/bandeja 10 kilos";          S = 1500; T = 10 }
    9  = @{ D = 44292; L = "Primera"; M = 50; N = 14000; O = 14000; P = 14000; Q = "`$/bandeja 10 kilos";          S = 1400; T = 10 }
    10 = @{ D = 44391; L = "Primera"; M = 15; N = 1500;  O = 1500;  P = 1500;  Q = "`$/kilo (en caja de 15 kilos)"; S = 1500; T = 1 }
    11 = @{ D = 44391; L = "Segunda"; M = 20; N = 1000;  O = 1000;  P = 1000;  Q = "`$/kilo (en caja de 15 kilos)"; S = 1000; T = 1 }
}

foreach ($r in $rows.Keys) {
    $data = $rows[$r]

    $ws.Cells.Item($r, 4).Value  = $data.D    # D: Fecha
    $ws.Cells.Item($r, 12).Value = $data.L    # L: Calidad
    $ws.Cells.Item($r, 13).Value = $data.M    # M: Volumen
    $ws.Cells.Item($r, 14).Value = $data.N    # N: Precio minimo
    $ws.Cells.Item($r, 15).Value = $data.O    # O: Precio maximo
    $ws.Cells.Item($r, 16).Value = $data.P    # P: Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $data.Q    # Q: Unidad de comercializacion
    $ws.Cells.Item($r, 19).Value = $data.S    # S: Precio $/Kg
    $ws.Cells.Item($r, 20).Value = $data.T    # T: Kg / unidad
}
